# Update cryptos list cell values per upstream diff
# (prices/volumes refreshed by the scheduled GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price column values are purely numeric-looking (e.g. "409.44"); the
# source sheet stores them as text, so force text formatting before writing
# to stop Excel's COM layer from auto-coercing them into numbers.
$textCells = @("D5","D6","D9","D11","D12","D13","D16","D18","D21","D22","D24","D25","D26","D27","D28","D31","D34","D37","D38","D41","D43","D44","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.984.85"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "3.412.09"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "409.44"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "129.37"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  +6.45%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.734"
$ws.Range("E9").Value = "  +5.73%  "
$ws.Range("E10").Value = "  +3.59%  "
$ws.Range("D11").Value = "42.80"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "9.25"
$ws.Range("E12").Value = "  +9.48%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  +42.64%  "
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "3.952.97"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "21.24"
$ws.Range("E16").Value = "  +7.01%  "
$ws.Range("D17").Value = "3.410.23"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "12.53"
$ws.Range("E18").Value = "  +8.28%  "
$ws.Range("E19").Value = "  +6.72%  "
$ws.Range("D20").Value = "61.907.25"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").Value = "453.27"
$ws.Range("E21").Value = "  +44.36%  "
$ws.Range("D22").Value = "92.05"
$ws.Range("E22").Value = "  +8.80%  "
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "13.26"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("D25").Value = "3.30"
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("D26").Value = "9.38"
$ws.Range("E26").Value = "  +14.91%  "
$ws.Range("D27").Value = "33.21"
$ws.Range("E27").Value = "  +11.42%  "
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").Value = "12.04"
$ws.Range("E31").Value = "  +5.65%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "42.87"
$ws.Range("E34").Value = "  -4.17%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("D37").Value = "53.81"
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  +7.75%  "
$ws.Range("D41").Value = "2.95"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Value = "142.75"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "4.25"
$ws.Range("E44").Value = "  +8.31%  "
$ws.Range("E45").Value = "  +15.94%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "16.65"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "0.149"
$ws.Range("E48").Value = "  +23.66%  "
$ws.Range("D49").Value = "22.43"
$ws.Range("E49").Value = "  +4.94%  "
$ws.Range("D50").Value = "2.14"
$ws.Range("E50").Value = "  +8.65%  "
$ws.Range("D51").Value = "3.756.17"
$ws.Range("E51").Value = "  -0.72%  "

# Restore the default style on the forced-text cells (only the number format
# was special-cased above; the workbook's original cells carry no explicit style).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

